# Apply updated cryptocurrency price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.646.25"
$ws.Range("E2").Value = "  -3.35%  "
$ws.Range("D3").Value = "2.604.31"
$ws.Range("E3").Value = "  -2.27%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'572.63"
$ws.Range("E5").Value = "  -4.28%  "
$ws.Range("D6").Value = "'155.45"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -6.93%  "
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("E11").Value = "  -5.09%  "
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").Value = "'28.16"
$ws.Range("E13").Value = "  -3.10%  "
$ws.Range("D14").Value = "3.074.13"
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("E15").Value = "  -8.32%  "
$ws.Range("D16").Value = "63.445.92"
$ws.Range("E16").Value = "  -3.53%  "
$ws.Range("D17").Value = "2.609.54"
$ws.Range("E17").Value = "  -2.02%  "
$ws.Range("E18").Value = "  -5.26%  "
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("E20").Value = "  -5.16%  "
$ws.Range("D21").Value = "'342.80"
$ws.Range("E21").Value = "  -2.59%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "'67.21"
$ws.Range("E23").Value = "  -3.77%  "
$ws.Range("E24").Value = "  -2.43%  "
$ws.Range("E25").Value = "  -3.23%  "
$ws.Range("D26").Value = "'587.91"
$ws.Range("E26").Value = "  +2.17%  "
$ws.Range("D27").Value = "'9.15"
$ws.Range("E27").Value = "  -6.12%  "
$ws.Range("E28").Value = "  -3.84%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "'0.161"
$ws.Range("E30").Value = "  -1.88%  "
$ws.Range("E31").Value = "  -3.90%  "
$ws.Range("E32").Value = "  -3.93%  "
$ws.Range("E33").Value = "  -5.34%  "
$ws.Range("E34").Value = "  -2.53%  "
$ws.Range("D35").Value = "'5.42"
$ws.Range("E35").Value = "  -3.00%  "
$ws.Range("D36").Value = "'0.405"
$ws.Range("E36").Value = "  -4.38%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'19.69"
$ws.Range("E38").Value = "  -4.20%  "
$ws.Range("D39").Value = "'155.38"
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("E40").Value = "  -4.57%  "
$ws.Range("D42").Value = "'2.51"
$ws.Range("E42").Value = "  +7.83%  "
$ws.Range("D43").Value = "'41.27"
$ws.Range("E43").Value = "  -3.46%  "
$ws.Range("D44").Value = "'156.29"
$ws.Range("E44").Value = "  -3.07%  "
$ws.Range("E45").Value = "  -4.50%  "
$ws.Range("D46").Value = "'23.12"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").Value = "'0.0588"
$ws.Range("E47").Value = "  -5.00%  "
$ws.Range("E48").Value = "  -2.44%  "
$ws.Range("E49").Value = "  -1.70%  "
$ws.Range("E50").Value = "  -3.89%  "
$ws.Range("D51").Value = "'18.78"
$ws.Range("E51").Value = "  -5.36%  "

# Clear the quote-prefix styling artifact so these cells keep their original (unstyled) look
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D51").Style = "Normal"

Write-Output "Updated cryptos list"
